$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.9693716918425304
$ws.Range("J2").Value = 0.9693716918425304
$ws.Range("M2").Value = 15.24491733333333
$ws.Range("N2").Value = 45.73475199999999
$ws.Range("O2").Value = 0.4831257321597052
$ws.Range("P2").Value = 0.4831257321597052
$ws.Range("Q2").Value = 139.5970880613618
$ws.Range("R2").Value = 1256.373792552256
$ws.Range("S2").Value = 0.4683284083563146
$ws.Range("T2").Value = 0.4683284083563146

# Row 3
$ws.Range("I3").Value = 0.9693716918425304
$ws.Range("J3").Value = 0.9693716918425304
$ws.Range("O3").Value = 0.327710667227878
$ws.Range("P3").Value = 0.327710667227878
$ws.Range("Q3").Value = 94.69057809683203
$ws.Range("S3").Value = 0.3176734439255326
$ws.Range("T3").Value = 0.3176734439255326

# Row 4
$ws.Range("I4").Value = 0.9693716918425304
$ws.Range("J4").Value = 0.9693716918425304
$ws.Range("M4").Value = 5.969012333333333
$ws.Range("N4").Value = 17.907037
$ws.Range("O4").Value = 0.1891636006124168
$ws.Range("P4").Value = 0.1891636006124168
$ws.Range("Q4").Value = 54.65800319649845
$ws.Range("R4").Value = 491.922028768486
$ws.Range("S4").Value = 0.1833698395606832
$ws.Range("T4").Value = 0.1833698395606832

# Row 5
$ws.Range("G5").Value = 0.2893236666666667
$ws.Range("H5").Value = 0.867971
$ws.Range("I5").Value = 0.03062830815746963
$ws.Range("J5").Value = 0.03062830815746962
$ws.Range("M5").Value = 15.24491733333333
$ws.Range("N5").Value = 45.73475199999999
$ws.Range("O5").Value = 0.4831257321597052
$ws.Range("P5").Value = 0.4831257321597052
$ws.Range("Q5").Value = 4.410715380910222
$ws.Range("R5").Value = 39.696438428192
$ws.Range("S5").Value = 0.01479732380339059
$ws.Range("T5").Value = 0.01479732380339058

# Row 6
$ws.Range("G6").Value = 0.2893236666666667
$ws.Range("H6").Value = 0.867971
$ws.Range("I6").Value = 0.03062830815746963
$ws.Range("J6").Value = 0.03062830815746962
$ws.Range("O6").Value = 0.327710667227878
$ws.Range("P6").Value = 0.327710667227878
$ws.Range("Q6").Value = 2.991847430624
$ws.Range("R6").Value = 26.926626875616
$ws.Range("S6").Value = 0.01003722330234543
$ws.Range("T6").Value = 0.01003722330234543

# Row 7
$ws.Range("G7").Value = 0.2893236666666667
$ws.Range("H7").Value = 0.867971
$ws.Range("I7").Value = 0.03062830815746963
$ws.Range("J7").Value = 0.03062830815746962
$ws.Range("M7").Value = 5.969012333333333
$ws.Range("N7").Value = 17.907037
$ws.Range("O7").Value = 0.1891636006124168
$ws.Range("P7").Value = 0.1891636006124168
$ws.Range("Q7").Value = 1.726976534658555
$ws.Range("R7").Value = 15.542788811927
$ws.Range("S7").Value = 0.005793761051733612
$ws.Range("T7").Value = 0.005793761051733612
